$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -0.7
$ws.Range("B4").Value = -0.7
$ws.Range("D9").Value = 1.5
$ws.Range("H9").Value = 2

$ws.Range("L17").Select()
